$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the placeholder "NA" value in C194 (becomes an empty cell, just like
# other rows that have no page number, e.g. C188).
$ws.Range("C194").Value = ""

# New row 195: 2025-08-04 / développement durable / page 124
$c = $ws.Range("A195")
$c.NumberFormat = "@"          # keep the date-looking text as literal text
$c.Value = "2025-08-04"
$c.Style = "Normal"            # drop the Text number format again, cell stays a string
$ws.Range("B195").Value = "développement durable"
$ws.Range("C195").Value = 124
$ws.Range("D195").Value = 1

# New row 196: 2025-08-04 / bonnes pratiques / page 126
$c = $ws.Range("A196")
$c.NumberFormat = "@"
$c.Value = "2025-08-04"
$c.Style = "Normal"
$ws.Range("B196").Value = "bonnes pratiques"
$ws.Range("C196").Value = 126
$ws.Range("D196").Value = 1
